$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("survey") ---
$ws1 = $wb.Worksheets.Item(1)

# Remove the obsolete header columns I:P (display.hint, choice_filter,
# hideInContents, display.audio, display.video, constraint,
# constraint_message, templatePath) - this also removes their data.
$ws1.Range("I1:P1").EntireColumn.Delete()

# Rename the remaining "display.text" header to "display.prompt.text"
$ws1.Range("H1").Value = "display.prompt.text"

# Resize the first three columns
$ws1.Range("A1").EntireColumn.ColumnWidth = 12.8
$ws1.Range("B1").EntireColumn.ColumnWidth = 9.8
$ws1.Range("C1").EntireColumn.ColumnWidth = 10.3

# --- Sheet 3 ("settings") ---
$ws3 = $wb.Worksheets.Item(3)

# Rename "display.title" header to "display.title.text"
$ws3.Range("C1").Value = "display.title.text"

# Give column A an explicit width
$ws3.Range("A1").EntireColumn.ColumnWidth = 23.8

# --- Update selections shown in each sheet view ---
# (select on sheet1 first, then re-select sheet3 so sheet3 stays the
# active/visible tab as in the original workbook)
$null = $ws1.Range("J11").Select()
$null = $ws3.Select()
$null = $ws3.Range("C2").Select()
